# feat: add 2022-Q4 data
#
# The workbook originally has two sheets: "总计" and "2022-Q3".
# This script:
#   1. Duplicates the existing "2022-Q3" sheet so the old data is preserved
#      in a new sheet (which stays named "2022-Q3").
#   2. Renames/repurposes the original sheet to "2022-Q4" and replaces its
#      contents with the new quarter's fund-holding data.
#   3. Updates the "总计" (totals) summary sheet with a new row for
#      2022-Q4, pushing the existing 2022-Q3 row down.

$wb = $excel.ActiveWorkbook

$total = $wb.Worksheets.Item(1)   # "总计"
$orig  = $wb.Worksheets.Item(2)   # currently "2022-Q3"

# --- 1. Duplicate the "2022-Q3" sheet, placing the copy right after it ---
$orig.Copy($null, $orig)

$q4 = $wb.Worksheets.Item(2)      # keeps the original sheet identity
$q3 = $wb.Worksheets.Item(3)      # the new copy, keeps the old data

$q4.Name = "2022-Q4"
$q3.Name = "2022-Q3"

# Restore the originally-selected tab onto the (content-wise unchanged)
# "2022-Q3" sheet.
$q3.Activate()

# --- 2. Replace the contents of the (renamed) 2022-Q4 sheet -------------
$q4.Cells.ClearContents()

$headers = @("基金代码","基金名称","基金规模","股票总仓位","仓位占比","持有市值(亿元)","仓位排名")
for ($c = 0; $c -lt $headers.Length; $c++) {
    $q4.Cells.Item(1, 2 + $c).Value = $headers[$c]
}

$rows = @(
    @("010405","惠升医药健康6个月持有期混合","9.98","72.67","2.34","0.2335",10),
    @("470888","汇添富香港优势精选混合（QDII）","2.68","93.08","7.87","0.2109",2),
    @("006603","嘉实互融精选股票","0.63","91.84","7.01","0.0442",3),
    @("005520","国投瑞银创新医疗混合","0.49","92.50","3.60","0.0176",7),
    @("006923","前海开源沪港深非周期性行业股票A","0.28","90.65","4.43","0.0124",9),
    @("009128","明亚价值长青混合A","0.39","52.21","2.80","0.0109",4),
    @("006924","前海开源沪港深非周期性行业股票C","0.24","90.65","4.43","0.0106",9),
    @("161124","易方达香港恒生综合小型股指数（QDII-LOF）A","0.24","94.45","1.25","0.0030",10),
    @("006263","易方达香港恒生综合小型股指数（QDII-LOF）C","0.05","94.45","1.25","0.0006",10),
    @("009129","明亚价值长青混合C","0.00","52.21","2.80",$null,4)
)

# Force columns B, D, E, F and (almost all of) G to store their numeric-
# looking values as genuine text, matching the source data, instead of
# letting them silently get auto-converted to numbers.
$q4.Range("B2:B11").NumberFormat = "@"
$q4.Range("D2:D11").NumberFormat = "@"
$q4.Range("E2:E11").NumberFormat = "@"
$q4.Range("F2:F11").NumberFormat = "@"
$q4.Range("G2:G10").NumberFormat = "@"

for ($i = 0; $i -lt $rows.Length; $i++) {
    $r = 2 + $i
    $row = $rows[$i]

    $q4.Cells.Item($r, 1).Value = $i
    $q4.Cells.Item($r, 2).Value = $row[0]
    $q4.Cells.Item($r, 3).Value = $row[1]
    $q4.Cells.Item($r, 4).Value = $row[2]
    $q4.Cells.Item($r, 5).Value = $row[3]
    $q4.Cells.Item($r, 6).Value = $row[4]

    if ($i -eq 9) {
        $q4.Cells.Item($r, 7).Value = 0
    } else {
        $q4.Cells.Item($r, 7).Value = $row[5]
    }

    $q4.Cells.Item($r, 8).Value = $row[6]
}

# Re-apply the sheet's normal (default) formatting on top of the values we
# just wrote, so the forced "@" text format above doesn't leave a stray
# number-format behind on the cells (it only influenced how the value was
# interpreted at write time).
$total.Range("C2").Copy()
$q4.Range("B2:H11").PasteSpecial(-4122)

# Header row (B1:H1) and the index column (A2:A11) use the same bold /
# bordered style as the other sheets.
$total.Range("B1").Copy()
$q4.Range("B1:H1").PasteSpecial(-4122)

$total.Range("A2").Copy()
$q4.Range("A2:A11").PasteSpecial(-4122)

# --- 3. Update the "总计" summary sheet ----------------------------------
# Give row 3 (the shifted-down 2022-Q3 row) the same index-column style as
# row 2 before writing the new values into both rows.
$total.Range("A2").Copy()
$total.Range("A3").PasteSpecial(-4122)

$total.Cells.Item(3, 1).Value = 1
$total.Cells.Item(3, 2).Value = "2022-Q3"
$total.Cells.Item(3, 3).Value = 4
$total.Cells.Item(3, 4).Value = 0.36

$total.Cells.Item(2, 1).Value = 0
$total.Cells.Item(2, 2).Value = "2022-Q4"
$total.Cells.Item(2, 3).Value = 10
$total.Cells.Item(2, 4).Value = 0.54
